$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Add the new "short_term_investment_minimal_irr" row under the existing Coupling Parameters
$ws.Range("A7").Value = "short_term_investment_minimal_irr"
$ws.Range("B7").Value = 0.3

# Widen column A to fit the new longer label
$ws.Columns.Item(1).ColumnWidth = 32

# Move the active selection to J8, matching the saved view state
$ws.Range("J8").Select() | Out-Null
